$d = $word.ActiveDocument
$t = $d.Tables.Item(2)

# Shift the Assignment/Points rows (rows 2-5) up by one position,
# i.e. each row takes on the assignment name + points of the row
# that follows it (row 5 -> row 2's new values come from row 3, etc.)
# New order: Multimedia Resume/20, Mini App/30, Web Application/40, Self Evaluation/10

$t.Cell(2, 1).Range.Text = "Multimedia Resume"
$t.Cell(2, 2).Range.Text = "20"

$t.Cell(3, 1).Range.Text = "Mini App"
$t.Cell(3, 2).Range.Text = "30"

$t.Cell(4, 1).Range.Text = "Web Application"
$t.Cell(4, 2).Range.Text = "40"

$t.Cell(5, 1).Range.Text = "Self Evaluation"
$t.Cell(5, 2).Range.Text = "10"
